$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "10-ago" column (AT) to the right of the existing data (AS).
$ws.Cells.Item(1, 46).Value = "10-ago"
$ws.Cells.Item(1, 46).NumberFormat = "@"

$newColumnValues = @(16, 15, 12, 14, 10, 15, 16, 22, 11, 20)
for ($i = 0; $i -lt $newColumnValues.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 46)
    $cell.Value = $newColumnValues[$i]
    $cell.HorizontalAlignment = -4108
    $cell.NumberFormat = "0"
}

# Remove the embedded picture from the sheet.
foreach ($shp in $ws.Shapes) {
    $shp.Delete()
}

# Match the author's final selection.
[void]$ws.Range("AT12").Select()
